$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

$ws.Range("A56").Value = 11
$ws.Range("B56").Value = 'Vega Monumental Concepción'
$ws.Range("C56").Value = 'Bíobío'
$ws.Range("D56").Value = 44757
$ws.Range("D56").NumberFormat = $ws.Range("D57").NumberFormat
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 100112032
$ws.Range("G56").Value = 'Zapallo italiano'
$ws.Range("H56").Value = 'Sin especificar'
$ws.Range("I56").Value = 'Primera'
$ws.Range("J56").Value = 100
$ws.Range("K56").Value = 9000
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = 9500
$ws.Range("N56").Value = '$/caja 50 unidades'
$ws.Range("O56").Value = 'Región de Arica y Parinacota'
$ws.Range("P56").Value = 190
$ws.Range("Q56").Value = 50
$ws.Range("R56").Value = 'Hortaliza'
